# Auto-generated Excel COM-interop script
# Applies updates to currentAveragePrice/LevePrice/LeveProfit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets,
# matching the upstream "Sheets via scheduled runner" data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43 (ALC)
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# Row 51 (ALC)
$ws.Range("H51").Value = 6500
$ws.Range("J51").Value = 6500
$ws.Range("L51").Value = 6500
$ws.Range("N51").Value = -7468

# Row 69 (ALC)
$ws.Range("H69").Value = 7345.6055
$ws.Range("J69").Value = 7345.6055
$ws.Range("L69").Value = 22036.8165
$ws.Range("N69").Value = -23784.8165

# Row 72 (ALC)
$ws.Range("H72").Value = 7345.6055
$ws.Range("J72").Value = 7345.6055
$ws.Range("L72").Value = 66110.4495
$ws.Range("N72").Value = -74846.4495

# Row 87 (ALC)
$ws.Range("H87").Value = 99999
$ws.Range("J87").Value = 99999
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102495

# Row 90 (ALC)
$ws.Range("H90").Value = 99999
$ws.Range("J90").Value = 99999
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477

$ws = $wb.Worksheets.Item("ARM")
# Row 37 (ARM)
$ws.Range("H37").Value = 43333
$ws.Range("I37").Value = 15000
$ws.Range("J37").Value = 99999
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 99999
$ws.Range("M37").Value = -14727
$ws.Range("N37").Value = -100545

# Row 63 (ARM)
$ws.Range("H63").Value = 1979.9
$ws.Range("I63").Value = 1987.375
$ws.Range("J63").Value = 1950
$ws.Range("K63").Value = 1987.375
$ws.Range("L63").Value = 1950
$ws.Range("M63").Value = -1301.375
$ws.Range("N63").Value = -3322

# Row 66 (ARM)
$ws.Range("H66").Value = 1979.9
$ws.Range("I66").Value = 1987.375
$ws.Range("J66").Value = 1950
$ws.Range("K66").Value = 9936.875
$ws.Range("L66").Value = 9750
$ws.Range("M66").Value = -6504.875
$ws.Range("N66").Value = -16614

# Row 74 (ARM)
$ws.Range("H74").Value = 4129.5
$ws.Range("I74").Value = 3688.889
$ws.Range("K74").Value = 3688.889
$ws.Range("M74").Value = -2814.889

# Row 77 (ARM)
$ws.Range("H77").Value = 4129.5
$ws.Range("I77").Value = 3688.889
$ws.Range("K77").Value = 18444.445
$ws.Range("M77").Value = -14076.445

$ws = $wb.Worksheets.Item("BSM")
# Row 26 (BSM)
$ws.Range("H26").Value = 45157
$ws.Range("I26").Value = 45157
$ws.Range("K26").Value = 45157
$ws.Range("M26").Value = -44865

# Row 96 (BSM)
$ws.Range("H96").Value = 21475.334
$ws.Range("I96").Value = 21770.4
$ws.Range("K96").Value = 21770.4
$ws.Range("M96").Value = -19024.4

# Row 99 (BSM)
$ws.Range("H99").Value = 166667660
$ws.Range("I99").Value = 166667660
$ws.Range("K99").Value = 166667660
$ws.Range("M99").Value = -166666162

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 2062.25
$ws.Range("I22").Value = 1174.5
$ws.Range("J22").Value = 2950
$ws.Range("K22").Value = 1174.5
$ws.Range("L22").Value = 2950
$ws.Range("M22").Value = -824.5
$ws.Range("N22").Value = -3650

# Row 55 (CRP)
$ws.Range("H55").Value = 15464.833
$ws.Range("I55").Value = 6400
$ws.Range("J55").Value = 19997.25
$ws.Range("K55").Value = 6400
$ws.Range("L55").Value = 19997.25
$ws.Range("M55").Value = -6085
$ws.Range("N55").Value = -20627.25

# Row 58 (CRP)
$ws.Range("H58").Value = 2514.6667
$ws.Range("I58").Value = 1704.625
$ws.Range("K58").Value = 1704.625
$ws.Range("M58").Value = -1501.625

# Row 107 (CRP)
$ws.Range("H107").Value = 455
$ws.Range("J107").Value = 525.5
$ws.Range("L107").Value = 525.5
$ws.Range("N107").Value = -4365.5

# Row 132 (CRP)
$ws.Range("H132").Value = 1711.6923
$ws.Range("I132").Value = 1437.7916
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 4313.3748
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -1783.3748
$ws.Range("N132").Value = -20055.5

# Row 134 (CRP)
$ws.Range("H134").Value = 1792.625
$ws.Range("I134").Value = 654.5
$ws.Range("J134").Value = 5207
$ws.Range("K134").Value = 1963.5
$ws.Range("L134").Value = 15621
$ws.Range("M134").Value = 571.5
$ws.Range("N134").Value = -20691

# Row 136 (CRP)
$ws.Range("H136").Value = 2514.6667
$ws.Range("I136").Value = 1704.625
$ws.Range("K136").Value = 5113.875
$ws.Range("M136").Value = -2563.875

$ws = $wb.Worksheets.Item("CUL")
# Row 57 (CUL)
$ws.Range("H57").Value = 1806.25
$ws.Range("I57").Value = 1408.3334
$ws.Range("K57").Value = 4225.0002
$ws.Range("M57").Value = -3666.0002

# Row 103 (CUL)
$ws.Range("H103").Value = 1165.1
$ws.Range("I103").Value = 627
$ws.Range("J103").Value = 1224.8889
$ws.Range("K103").Value = 1881
$ws.Range("L103").Value = 3674.6667
$ws.Range("M103").Value = -1002
$ws.Range("N103").Value = -5432.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 4386.143
$ws.Range("I80").Value = 4150.75
$ws.Range("J80").Value = 4700
$ws.Range("K80").Value = 4150.75
$ws.Range("L80").Value = 4700
$ws.Range("M80").Value = -3152.75
$ws.Range("N80").Value = -6696

# Row 83 (GSM)
$ws.Range("H83").Value = 4386.143
$ws.Range("I83").Value = 4150.75
$ws.Range("J83").Value = 4700
$ws.Range("K83").Value = 20753.75
$ws.Range("L83").Value = 23500
$ws.Range("M83").Value = -15761.75
$ws.Range("N83").Value = -33484

# Row 99 (GSM)
$ws.Range("H99").Value = 4600
$ws.Range("I99").Value = 4600
$ws.Range("K99").Value = 4600
$ws.Range("M99").Value = -2354

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 2238
$ws.Range("J22").Value = 2875
$ws.Range("L22").Value = 2875
$ws.Range("N22").Value = -3465

# Row 27 (LTW)
$ws.Range("H27").Value = 2238
$ws.Range("J27").Value = 2875
$ws.Range("L27").Value = 2875
$ws.Range("N27").Value = -3089

# Row 68 (LTW)
$ws.Range("H68").Value = 7704
$ws.Range("J68").Value = 9106.857
$ws.Range("L68").Value = 9106.857
$ws.Range("N68").Value = -10604.857

# Row 71 (LTW)
$ws.Range("H71").Value = 7704
$ws.Range("J71").Value = 9106.857
$ws.Range("L71").Value = 45534.285
$ws.Range("N71").Value = -53022.285

# Row 93 (LTW)
$ws.Range("H93").Value = 1328.625
$ws.Range("I93").Value = 971.5
$ws.Range("K93").Value = 971.5
$ws.Range("M93").Value = 276.5

# Row 136 (LTW)
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# Row 140 (LTW)
$ws.Range("H140").Value = 137500
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
# Row 18 (WVR)
$ws.Range("H18").Value = 16400
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 16400
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 16400
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -16746

# Row 81 (WVR)
$ws.Range("H81").Value = 956.2
$ws.Range("I81").Value = 956.2
$ws.Range("K81").Value = 1912.4
$ws.Range("M81").Value = -851.4000000000001

# Row 84 (WVR)
$ws.Range("H84").Value = 956.2
$ws.Range("I84").Value = 956.2
$ws.Range("K84").Value = 9562
$ws.Range("M84").Value = -4258

# Row 132 (WVR)
$ws.Range("H132").Value = 6249.5
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 7499
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 22497
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -27557

# Row 136 (WVR)
$ws.Range("H136").Value = 2279.7856
$ws.Range("I136").Value = 1596.3823
$ws.Range("K136").Value = 4789.1469
$ws.Range("M136").Value = -2239.1469
